# Update gh-pages to output generated at 456a3b4
# Apply the updated "想去人数" (F column) figures to both the "展览"
# and "全部类型" worksheets (which hold duplicate data).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F7").Value = 11959
    $ws.Range("F11").Value = 411
    $ws.Range("F15").Value = 13445
    $ws.Range("F20").Value = 287
    $ws.Range("F23").Value = 168
}
